$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current values for rows 2-14, columns D,K,L,M,N,O,P,Q,R,S
$r2_D = $ws.Cells.Item(2, 4).Value2()
$r2_K = $ws.Cells.Item(2, 11).Value2()
$r2_L = $ws.Cells.Item(2, 12).Value2()
$r2_M = $ws.Cells.Item(2, 13).Value2()
$r2_N = $ws.Cells.Item(2, 14).Value2()
$r2_O = $ws.Cells.Item(2, 15).Value2()
$r2_P = $ws.Cells.Item(2, 16).Value2()
$r2_Q = $ws.Cells.Item(2, 17).Value2()
$r2_R = $ws.Cells.Item(2, 18).Value2()
$r2_S = $ws.Cells.Item(2, 19).Value2()

$r3_D = $ws.Cells.Item(3, 4).Value2()
$r3_K = $ws.Cells.Item(3, 11).Value2()
$r3_L = $ws.Cells.Item(3, 12).Value2()
$r3_M = $ws.Cells.Item(3, 13).Value2()
$r3_N = $ws.Cells.Item(3, 14).Value2()
$r3_O = $ws.Cells.Item(3, 15).Value2()
$r3_P = $ws.Cells.Item(3, 16).Value2()
$r3_Q = $ws.Cells.Item(3, 17).Value2()
$r3_R = $ws.Cells.Item(3, 18).Value2()
$r3_S = $ws.Cells.Item(3, 19).Value2()

$r4_D = $ws.Cells.Item(4, 4).Value2()
$r4_K = $ws.Cells.Item(4, 11).Value2()
$r4_L = $ws.Cells.Item(4, 12).Value2()
$r4_M = $ws.Cells.Item(4, 13).Value2()
$r4_N = $ws.Cells.Item(4, 14).Value2()
$r4_O = $ws.Cells.Item(4, 15).Value2()
$r4_P = $ws.Cells.Item(4, 16).Value2()
$r4_Q = $ws.Cells.Item(4, 17).Value2()
$r4_R = $ws.Cells.Item(4, 18).Value2()
$r4_S = $ws.Cells.Item(4, 19).Value2()

$r5_D = $ws.Cells.Item(5, 4).Value2()
$r5_K = $ws.Cells.Item(5, 11).Value2()
$r5_L = $ws.Cells.Item(5, 12).Value2()
$r5_M = $ws.Cells.Item(5, 13).Value2()
$r5_N = $ws.Cells.Item(5, 14).Value2()
$r5_O = $ws.Cells.Item(5, 15).Value2()
$r5_P = $ws.Cells.Item(5, 16).Value2()
$r5_Q = $ws.Cells.Item(5, 17).Value2()
$r5_R = $ws.Cells.Item(5, 18).Value2()
$r5_S = $ws.Cells.Item(5, 19).Value2()

$r6_D = $ws.Cells.Item(6, 4).Value2()
$r6_K = $ws.Cells.Item(6, 11).Value2()
$r6_L = $ws.Cells.Item(6, 12).Value2()
$r6_M = $ws.Cells.Item(6, 13).Value2()
$r6_N = $ws.Cells.Item(6, 14).Value2()
$r6_O = $ws.Cells.Item(6, 15).Value2()
$r6_P = $ws.Cells.Item(6, 16).Value2()
$r6_Q = $ws.Cells.Item(6, 17).Value2()
$r6_R = $ws.Cells.Item(6, 18).Value2()
$r6_S = $ws.Cells.Item(6, 19).Value2()

$r7_D = $ws.Cells.Item(7, 4).Value2()
$r7_K = $ws.Cells.Item(7, 11).Value2()
$r7_L = $ws.Cells.Item(7, 12).Value2()
$r7_M = $ws.Cells.Item(7, 13).Value2()
$r7_N = $ws.Cells.Item(7, 14).Value2()
$r7_O = $ws.Cells.Item(7, 15).Value2()
$r7_P = $ws.Cells.Item(7, 16).Value2()
$r7_Q = $ws.Cells.Item(7, 17).Value2()
$r7_R = $ws.Cells.Item(7, 18).Value2()
$r7_S = $ws.Cells.Item(7, 19).Value2()

$r8_D = $ws.Cells.Item(8, 4).Value2()
$r8_K = $ws.Cells.Item(8, 11).Value2()
$r8_L = $ws.Cells.Item(8, 12).Value2()
$r8_M = $ws.Cells.Item(8, 13).Value2()
$r8_N = $ws.Cells.Item(8, 14).Value2()
$r8_O = $ws.Cells.Item(8, 15).Value2()
$r8_P = $ws.Cells.Item(8, 16).Value2()
$r8_Q = $ws.Cells.Item(8, 17).Value2()
$r8_R = $ws.Cells.Item(8, 18).Value2()
$r8_S = $ws.Cells.Item(8, 19).Value2()

$r9_D = $ws.Cells.Item(9, 4).Value2()
$r9_K = $ws.Cells.Item(9, 11).Value2()
$r9_L = $ws.Cells.Item(9, 12).Value2()
$r9_M = $ws.Cells.Item(9, 13).Value2()
$r9_N = $ws.Cells.Item(9, 14).Value2()
$r9_O = $ws.Cells.Item(9, 15).Value2()
$r9_P = $ws.Cells.Item(9, 16).Value2()
$r9_Q = $ws.Cells.Item(9, 17).Value2()
$r9_R = $ws.Cells.Item(9, 18).Value2()
$r9_S = $ws.Cells.Item(9, 19).Value2()

$r10_D = $ws.Cells.Item(10, 4).Value2()
$r10_K = $ws.Cells.Item(10, 11).Value2()
$r10_L = $ws.Cells.Item(10, 12).Value2()
$r10_M = $ws.Cells.Item(10, 13).Value2()
$r10_N = $ws.Cells.Item(10, 14).Value2()
$r10_O = $ws.Cells.Item(10, 15).Value2()
$r10_P = $ws.Cells.Item(10, 16).Value2()
$r10_Q = $ws.Cells.Item(10, 17).Value2()
$r10_R = $ws.Cells.Item(10, 18).Value2()
$r10_S = $ws.Cells.Item(10, 19).Value2()

$r11_D = $ws.Cells.Item(11, 4).Value2()
$r11_K = $ws.Cells.Item(11, 11).Value2()
$r11_L = $ws.Cells.Item(11, 12).Value2()
$r11_M = $ws.Cells.Item(11, 13).Value2()
$r11_N = $ws.Cells.Item(11, 14).Value2()
$r11_O = $ws.Cells.Item(11, 15).Value2()
$r11_P = $ws.Cells.Item(11, 16).Value2()
$r11_Q = $ws.Cells.Item(11, 17).Value2()
$r11_R = $ws.Cells.Item(11, 18).Value2()
$r11_S = $ws.Cells.Item(11, 19).Value2()

$r12_D = $ws.Cells.Item(12, 4).Value2()
$r12_K = $ws.Cells.Item(12, 11).Value2()
$r12_L = $ws.Cells.Item(12, 12).Value2()
$r12_M = $ws.Cells.Item(12, 13).Value2()
$r12_N = $ws.Cells.Item(12, 14).Value2()
$r12_O = $ws.Cells.Item(12, 15).Value2()
$r12_P = $ws.Cells.Item(12, 16).Value2()
$r12_Q = $ws.Cells.Item(12, 17).Value2()
$r12_R = $ws.Cells.Item(12, 18).Value2()
$r12_S = $ws.Cells.Item(12, 19).Value2()

$r13_D = $ws.Cells.Item(13, 4).Value2()
$r13_K = $ws.Cells.Item(13, 11).Value2()
$r13_L = $ws.Cells.Item(13, 12).Value2()
$r13_M = $ws.Cells.Item(13, 13).Value2()
$r13_N = $ws.Cells.Item(13, 14).Value2()
$r13_O = $ws.Cells.Item(13, 15).Value2()
$r13_P = $ws.Cells.Item(13, 16).Value2()
$r13_Q = $ws.Cells.Item(13, 17).Value2()
$r13_R = $ws.Cells.Item(13, 18).Value2()
$r13_S = $ws.Cells.Item(13, 19).Value2()

$r14_D = $ws.Cells.Item(14, 4).Value2()
$r14_K = $ws.Cells.Item(14, 11).Value2()
$r14_L = $ws.Cells.Item(14, 12).Value2()
$r14_M = $ws.Cells.Item(14, 13).Value2()
$r14_N = $ws.Cells.Item(14, 14).Value2()
$r14_O = $ws.Cells.Item(14, 15).Value2()
$r14_P = $ws.Cells.Item(14, 16).Value2()
$r14_Q = $ws.Cells.Item(14, 17).Value2()
$r14_R = $ws.Cells.Item(14, 18).Value2()
$r14_S = $ws.Cells.Item(14, 19).Value2()

# Write shuffled values back according to the row permutation mapping
$ws.Cells.Item(2, 4).Value = $r14_D
$ws.Cells.Item(2, 11).Value = $r14_K
$ws.Cells.Item(2, 12).Value = $r14_L
$ws.Cells.Item(2, 13).Value = $r14_M
$ws.Cells.Item(2, 14).Value = $r14_N
$ws.Cells.Item(2, 15).Value = $r14_O
$ws.Cells.Item(2, 16).Value = $r14_P
$ws.Cells.Item(2, 17).Value = $r14_Q
$ws.Cells.Item(2, 18).Value = $r14_R
$ws.Cells.Item(2, 19).Value = $r14_S

$ws.Cells.Item(3, 4).Value = $r10_D
$ws.Cells.Item(3, 11).Value = $r10_K
$ws.Cells.Item(3, 12).Value = $r10_L
$ws.Cells.Item(3, 13).Value = $r10_M
$ws.Cells.Item(3, 14).Value = $r10_N
$ws.Cells.Item(3, 15).Value = $r10_O
$ws.Cells.Item(3, 16).Value = $r10_P
$ws.Cells.Item(3, 17).Value = $r10_Q
$ws.Cells.Item(3, 18).Value = $r10_R
$ws.Cells.Item(3, 19).Value = $r10_S

$ws.Cells.Item(4, 4).Value = $r11_D
$ws.Cells.Item(4, 11).Value = $r11_K
$ws.Cells.Item(4, 12).Value = $r11_L
$ws.Cells.Item(4, 13).Value = $r11_M
$ws.Cells.Item(4, 14).Value = $r11_N
$ws.Cells.Item(4, 15).Value = $r11_O
$ws.Cells.Item(4, 16).Value = $r11_P
$ws.Cells.Item(4, 17).Value = $r11_Q
$ws.Cells.Item(4, 18).Value = $r11_R
$ws.Cells.Item(4, 19).Value = $r11_S

$ws.Cells.Item(5, 4).Value = $r6_D
$ws.Cells.Item(5, 11).Value = $r6_K
$ws.Cells.Item(5, 12).Value = $r6_L
$ws.Cells.Item(5, 13).Value = $r6_M
$ws.Cells.Item(5, 14).Value = $r6_N
$ws.Cells.Item(5, 15).Value = $r6_O
$ws.Cells.Item(5, 16).Value = $r6_P
$ws.Cells.Item(5, 17).Value = $r6_Q
$ws.Cells.Item(5, 18).Value = $r6_R
$ws.Cells.Item(5, 19).Value = $r6_S

$ws.Cells.Item(6, 4).Value = $r7_D
$ws.Cells.Item(6, 11).Value = $r7_K
$ws.Cells.Item(6, 12).Value = $r7_L
$ws.Cells.Item(6, 13).Value = $r7_M
$ws.Cells.Item(6, 14).Value = $r7_N
$ws.Cells.Item(6, 15).Value = $r7_O
$ws.Cells.Item(6, 16).Value = $r7_P
$ws.Cells.Item(6, 17).Value = $r7_Q
$ws.Cells.Item(6, 18).Value = $r7_R
$ws.Cells.Item(6, 19).Value = $r7_S

$ws.Cells.Item(7, 4).Value = $r13_D
$ws.Cells.Item(7, 11).Value = $r13_K
$ws.Cells.Item(7, 12).Value = $r13_L
$ws.Cells.Item(7, 13).Value = $r13_M
$ws.Cells.Item(7, 14).Value = $r13_N
$ws.Cells.Item(7, 15).Value = $r13_O
$ws.Cells.Item(7, 16).Value = $r13_P
$ws.Cells.Item(7, 17).Value = $r13_Q
$ws.Cells.Item(7, 18).Value = $r13_R
$ws.Cells.Item(7, 19).Value = $r13_S

$ws.Cells.Item(8, 4).Value = $r2_D
$ws.Cells.Item(8, 11).Value = $r2_K
$ws.Cells.Item(8, 12).Value = $r2_L
$ws.Cells.Item(8, 13).Value = $r2_M
$ws.Cells.Item(8, 14).Value = $r2_N
$ws.Cells.Item(8, 15).Value = $r2_O
$ws.Cells.Item(8, 16).Value = $r2_P
$ws.Cells.Item(8, 17).Value = $r2_Q
$ws.Cells.Item(8, 18).Value = $r2_R
$ws.Cells.Item(8, 19).Value = $r2_S

$ws.Cells.Item(9, 4).Value = $r9_D
$ws.Cells.Item(9, 11).Value = $r9_K
$ws.Cells.Item(9, 12).Value = $r9_L
$ws.Cells.Item(9, 13).Value = $r9_M
$ws.Cells.Item(9, 14).Value = $r9_N
$ws.Cells.Item(9, 15).Value = $r9_O
$ws.Cells.Item(9, 16).Value = $r9_P
$ws.Cells.Item(9, 17).Value = $r9_Q
$ws.Cells.Item(9, 18).Value = $r9_R
$ws.Cells.Item(9, 19).Value = $r9_S

$ws.Cells.Item(10, 4).Value = $r3_D
$ws.Cells.Item(10, 11).Value = $r3_K
$ws.Cells.Item(10, 12).Value = $r3_L
$ws.Cells.Item(10, 13).Value = $r3_M
$ws.Cells.Item(10, 14).Value = $r3_N
$ws.Cells.Item(10, 15).Value = $r3_O
$ws.Cells.Item(10, 16).Value = $r3_P
$ws.Cells.Item(10, 17).Value = $r3_Q
$ws.Cells.Item(10, 18).Value = $r3_R
$ws.Cells.Item(10, 19).Value = $r3_S

$ws.Cells.Item(11, 4).Value = $r5_D
$ws.Cells.Item(11, 11).Value = $r5_K
$ws.Cells.Item(11, 12).Value = $r5_L
$ws.Cells.Item(11, 13).Value = $r5_M
$ws.Cells.Item(11, 14).Value = $r5_N
$ws.Cells.Item(11, 15).Value = $r5_O
$ws.Cells.Item(11, 16).Value = $r5_P
$ws.Cells.Item(11, 17).Value = $r5_Q
$ws.Cells.Item(11, 18).Value = $r5_R
$ws.Cells.Item(11, 19).Value = $r5_S

$ws.Cells.Item(12, 4).Value = $r8_D
$ws.Cells.Item(12, 11).Value = $r8_K
$ws.Cells.Item(12, 12).Value = $r8_L
$ws.Cells.Item(12, 13).Value = $r8_M
$ws.Cells.Item(12, 14).Value = $r8_N
$ws.Cells.Item(12, 15).Value = $r8_O
$ws.Cells.Item(12, 16).Value = $r8_P
$ws.Cells.Item(12, 17).Value = $r8_Q
$ws.Cells.Item(12, 18).Value = $r8_R
$ws.Cells.Item(12, 19).Value = $r8_S

$ws.Cells.Item(13, 4).Value = $r12_D
$ws.Cells.Item(13, 11).Value = $r12_K
$ws.Cells.Item(13, 12).Value = $r12_L
$ws.Cells.Item(13, 13).Value = $r12_M
$ws.Cells.Item(13, 14).Value = $r12_N
$ws.Cells.Item(13, 15).Value = $r12_O
$ws.Cells.Item(13, 16).Value = $r12_P
$ws.Cells.Item(13, 17).Value = $r12_Q
$ws.Cells.Item(13, 18).Value = $r12_R
$ws.Cells.Item(13, 19).Value = $r12_S

$ws.Cells.Item(14, 4).Value = $r4_D
$ws.Cells.Item(14, 11).Value = $r4_K
$ws.Cells.Item(14, 12).Value = $r4_L
$ws.Cells.Item(14, 13).Value = $r4_M
$ws.Cells.Item(14, 14).Value = $r4_N
$ws.Cells.Item(14, 15).Value = $r4_O
$ws.Cells.Item(14, 16).Value = $r4_P
$ws.Cells.Item(14, 17).Value = $r4_Q
$ws.Cells.Item(14, 18).Value = $r4_R
$ws.Cells.Item(14, 19).Value = $r4_S
